$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the k column (J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Rows 14-17: summary labels + aggregate formulas
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the new summary values: bold, size-12 font, vertically centered.
# Apply directly to B14 (single fresh cell -> single clean style index),
# then propagate the same format to B15:B17 via copy/paste-special so we
# don't mint extra intermediate cellXfs entries.
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108
$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row height for the new label rows
$ws.Range("A14:B17").RowHeight = 15.6

# Print setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection left on the new summary block
$ws.Range("A14:B17").Select()
